# Code clean up, new way to create characters
#
# Rebuild the data table on Sheet1 using the new set of shared strings
# (Attention Deployment / Cognitive Change / [Situation Selection, Weakly] /
# [Attention Deployment, Strongly] / [Cognitive Change, Strongly] /
# [Response Modulation, Lightly]) and refresh the sampled A/C numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows 2-11 (MOOD / EMOTION / INTENSITY / EVENT / APPLIED STRATEGY) ---
$ws.Range("A2").Value = -0.7413855195045471
$ws.Range("B2").Value = "Distress"
$ws.Range("C2").Value = 2.387104034423828
$ws.Range("D2").Value = "Talk"
$ws.Range("E2").Value = "Attention Deployment"

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "Love"
$ws.Range("C3").Value = 1.7927955389022827
$ws.Range("D3").Value = "Hello"
$ws.Range("E3").Value = "None"

$ws.Range("A4").Value = 0.8303518295288086
$ws.Range("B4").Value = "Love"
$ws.Range("C4").Value = 2.6735565662384033
$ws.Range("D4").Value = "Conversation"
$ws.Range("E4").Value = "None"

$ws.Range("A5").Value = 2.3478033542633057
$ws.Range("B5").Value = "Love"
$ws.Range("C5").Value = 4.9165802001953125
$ws.Range("D5").Value = "Hug"
$ws.Range("E5").Value = "None"

$ws.Range("A6").Value = 1.046940803527832
$ws.Range("B6").Value = "Distress"
$ws.Range("C6").Value = 4.1016740798950195
$ws.Range("D6").Value = "Discussion"
$ws.Range("E6").Value = "None"

$ws.Range("A7").Value = 1.7211663722991943
$ws.Range("B7").Value = "Joy"
$ws.Range("C7").Value = 2.2095818519592285
$ws.Range("D7").Value = "Congrat"
$ws.Range("E7").Value = "None"

$ws.Range("A8").Value = 0.7149765491485596
$ws.Range("B8").Value = "Distress"
$ws.Range("C8").Value = 3.1760647296905518
$ws.Range("D8").Value = "Bye"
$ws.Range("E8").Value = "Attention Deployment"

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = "Hate"
$ws.Range("C9").Value = 2.5324015617370605
$ws.Range("D9").Value = "Fired"
$ws.Range("E9").Value = "Cognitive Change"

$ws.Range("A10").Value = -0.5182909965515137
$ws.Range("B10").Value = "Hate"
$ws.Range("C10").Value = 1.6687870025634766
$ws.Range("D10").Value = "Crash"
$ws.Range("E10").Value = "Cognitive Change"

$ws.Range("A11").Value = 1.8139854669570923
$ws.Range("B11").Value = "Joy"
$ws.Range("C11").Value = 7.490267276763916
$ws.Range("D11").Value = "Profits"
$ws.Range("E11").Value = "None"

# Rows 12-13 used to carry a full A-E record each; now they only hold a
# single PERSONALITY TRAITS entry in column F (see below), so drop the old
# A-E content that used to live here.
$ws.Range("A12:E13").ClearContents()

# --- PERSONALITY TRAITS list, column F, rows 12-16 ---
$ws.Range("F12").Value = "Low Conscientiousness"
$ws.Range("F13").Value = "Low Extraversion"
$ws.Range("F14").Value = "Low Neuroticism"
$ws.Range("F15").Value = "High Agreeableness"
$ws.Range("F16").Value = "Low Openness"

# --- STRATEGIES RELATED list moves from column F to column G, gains entries ---
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("G17").Value = "[Situation Selection, Weakly]"
$ws.Range("G18").Value = "[Situation Modification, Weakly]"
$ws.Range("G19").Value = "[Attention Deployment, Strongly]"
$ws.Range("G20").Value = "[Cognitive Change, Strongly]"
$ws.Range("G21").Value = "[Response Modulation, Lightly]"

# --- DOMINANT PERSONALITY moves from H20 to H22 ---
$ws.Range("H20").ClearContents()
$ws.Range("H22").Value = "Agreeableness"
